# C5-PowerPoint.pptx edit:
#   1. Slide 6's table switches to a different table style
#      ({87E43637-6EB6-4B04-B9E7-E0DDA21E4DE2} -> {86032DB1-8DF3-44B5-BB0B-A4274DFE1E3B}).
#   2. The deck's theme colour scheme (bound to the slide master / theme1.xml)
#      is swapped from the "Integral" palette to the stock "Office" palette.

function HexToComRGB([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# --- 1. Re-style the table on slide 6 -------------------------------------
$slide6 = $p.Slides.Item(6)
for ($i = 1; $i -le $slide6.Shapes.Count; $i++) {
    $shape = $slide6.Shapes.Item($i)
    if ($shape.HasTable) {
        $shape.Table.ApplyStyle("{86032DB1-8DF3-44B5-BB0B-A4274DFE1E3B}")
    }
}

# --- 2. Swap the theme colour scheme to the stock Office palette ----------
$anySlide = $p.Slides.Item(1)
$themeColors = $anySlide.ThemeColorScheme

$themeColors.Item(1).RGB  = HexToComRGB "000000"   # dk1
$themeColors.Item(2).RGB  = HexToComRGB "FFFFFF"   # lt1
$themeColors.Item(3).RGB  = HexToComRGB "44546A"   # dk2
$themeColors.Item(4).RGB  = HexToComRGB "E7E6E6"   # lt2
$themeColors.Item(5).RGB  = HexToComRGB "5B9BD5"   # accent1
$themeColors.Item(6).RGB  = HexToComRGB "ED7D31"   # accent2
$themeColors.Item(7).RGB  = HexToComRGB "A5A5A5"   # accent3
$themeColors.Item(8).RGB  = HexToComRGB "FFC000"   # accent4
$themeColors.Item(9).RGB  = HexToComRGB "4472C4"   # accent5
$themeColors.Item(10).RGB = HexToComRGB "70AD47"   # accent6
$themeColors.Item(11).RGB = HexToComRGB "0563C1"   # hyperlink
$themeColors.Item(12).RGB = HexToComRGB "954F72"   # followed hyperlink
